# LOM3260.xlsx update — "Build site at 2023-01-09 16:18:13 UTC"
#
#  - Ativação date bumped 01/01/2020 -> 01/01/2023              (rows 8 & 13)
#  - Responsible professor changed to Emerson Gonçalves de Melo (rows 10 & 15)
#  - English "Objectives:"    row (11) filled in
#  - English "Short syllabus:" row (14) filled in
#  - English "Syllabus:"       row (16) filled in
#  - "Norma de recuperação:" text expanded                      (row 20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a plain text value into a brand-new (previously empty) cell, copying
# the B/C number-format-free style (wrap text, top aligned; red font for the
# "modified" column C) from the row-3 reference cells so the generated
# cellXfs line up with the ones Excel already uses in this sheet.
function Set-NewCellText($addr, $fmtSrc, $value) {
    $ws.Range($fmtSrc).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $value
}

# Write a value that *looks* like a date/number (e.g. "01/01/2023") into a
# set of cells while keeping it stored as plain text, then restore each
# cell's original (numberformat-free) style so it keeps matching the B/C
# pair style. Cells sharing one format source are batched together so the
# transient "@" number format only mints one throwaway style per group.
function Set-TextValues($addrs, $fmtSrc, $value) {
    foreach ($addr in $addrs) {
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $value
    }
    foreach ($addr in $addrs) {
        $ws.Range($fmtSrc).Copy() | Out-Null
        $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    }
}

# --- Ativação: 01/01/2020 -> 01/01/2023 -------------------------------------
# (row 13 "Programa resumido:" mirrors the same date value in the source
# sheet, so it is updated in lock-step with row 8 "Ativação:")
Set-TextValues @("B8", "B13") "B3" "01/01/2023"
Set-TextValues @("C8", "C13") "C3" "01/01/2023"

# --- Objetivos: responsible professor ---------------------------------------
# (row 15 "Programa:" mirrors the same professor value in the source sheet)
$prof = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("B10").Value = $prof
$ws.Range("C10").Value = $prof
$ws.Range("B15").Value = $prof
$ws.Range("C15").Value = $prof

# --- Objectives: (English objectives text, row 11, was blank) --------------
$objectivesEn = "Provide the student with an introduction to modern scientific computing, using the Python language and its most popular numerical and graphical libraries: numpy, scipy, matplotlib, and pandas. At the end of the course, the student will be able to develop complex, small and medium-sized programs to solve engineering problems that involve numerical processing of large data sets and correlate variables using numerical methods."
Set-NewCellText "B11" "B3" $objectivesEn
Set-NewCellText "C11" "C3" $objectivesEn

# --- Short syllabus: (English short syllabus text, row 14, was blank) ------
$shortSyllabusEn = "Introduction to Python programming; keywords in Python; routines and functions; classes; numpy and the concept of slicing and indexing arrays; review of numerical methods using scipy; generating graphics and animations with the matplotlib library; creating graphical user interfaces using matplotlib.widgets"
Set-NewCellText "B14" "B3" $shortSyllabusEn
Set-NewCellText "C14" "C3" $shortSyllabusEn

# --- Syllabus: (English full syllabus text, row 16, was blank) -------------
$syllabusEn = "• Introduction to Python programming • Installing a Python distribution on Windows and Linux • Python file formatting • Conditional structures • Command loops • Other keywords and methods • Routines and functions • Multi-source code and personal libraries • Numerical and graphical libraries: numpy, scipy and matplotlib • Object-oriented programming: classes • Concept of objects and instances • Classes and subclasses • “Arrays” in numpy • The concept of array in numpy • “Slicing” and indexing • Working with files ( input and output) • Graphs in matplotlib • The matplotlib.pyplot library and 2D and 3D graphs • The matplotlib.animation library for creating animated graphs. • Graphical User Interface (GUI) • Simple interfaces with the matplotlib.widgets library."
Set-NewCellText "B16" "B3" $syllabusEn
Set-NewCellText "C16" "C3" $syllabusEn

# --- Norma de recuperação: expanded wording ---------------------------------
$recoveryNorm = "Média aritmética de exercícios e trabalhos propostos ao longo do curso e uma apresentação final de projeto."
$ws.Range("B20").Value = $recoveryNorm
$ws.Range("C20").Value = $recoveryNorm
